# Update column F (dSF) values to match repulled data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 5
$ws.Range("F6").Value = 0
$ws.Range("F9").Value = -8
$ws.Range("F11").Value = -5
$ws.Range("F14").Value = 1
$ws.Range("F16").Value = 3
